$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare header formatting for new columns K:O by copying format from existing header cell J1 ---
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:O1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Update header row (row 1) ---
$ws.Cells.Item(1, 8).Value = "Avg_Experiment_Time"   # H1: Experiment_Time -> Avg_Experiment_Time
$ws.Cells.Item(1, 9).Value = "Std_Total_Rounds"       # I1: Obs_Prob -> Std_Total_Rounds
$ws.Cells.Item(1, 10).Value = "Std_Expl_Cost"         # J1: Std_Total_Rounds -> Std_Expl_Cost
$ws.Cells.Item(1, 11).Value = "Std_Expl_Eff"          # K1: new
$ws.Cells.Item(1, 12).Value = "Std_Round_Time"        # L1: new
$ws.Cells.Item(1, 13).Value = "Std_Agent_Step_Time"   # M1: new
$ws.Cells.Item(1, 14).Value = "Std_Experiment_Time"   # N1: new
$ws.Cells.Item(1, 15).Value = "Obs_Prob"              # O1: new (moved from old I1)

# --- Update data rows 2-13 ---
# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 57.204
$ws.Cells.Item(2, 4).Value = 57.204
$ws.Cells.Item(2, 5).Value = 2.98872616
$ws.Cells.Item(2, 6).Value = 0.118009
$ws.Cells.Item(2, 7).Value = 0.118009
$ws.Cells.Item(2, 8).Value = 6.67244584
$ws.Cells.Item(2, 9).Value = 6.207047516956472
$ws.Cells.Item(2, 10).Value = 6.207047516956472
$ws.Cells.Item(2, 11).Value = 0.3197360522861278
$ws.Cells.Item(2, 12).Value = 0.01547734405374242
$ws.Cells.Item(2, 13).Value = 0.01547734405374242
$ws.Cells.Item(2, 14).Value = 0.4949738000287735
$ws.Cells.Item(2, 15).Value = 0.15

# Row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 90.636
$ws.Cells.Item(3, 4).Value = 90.636
$ws.Cells.Item(3, 5).Value = 1.89906912
$ws.Cells.Item(3, 6).Value = 0.07723984
$ws.Cells.Item(3, 7).Value = 0.07723984
$ws.Cells.Item(3, 8).Value = 6.89944312
$ws.Cells.Item(3, 9).Value = 12.2333484225397
$ws.Cells.Item(3, 10).Value = 12.2333484225397
$ws.Cells.Item(3, 11).Value = 0.2605441289568686
$ws.Cells.Item(3, 12).Value = 0.01118271574517007
$ws.Cells.Item(3, 13).Value = 0.01118271574517007
$ws.Cells.Item(3, 14).Value = 0.6730569641552351
$ws.Cells.Item(3, 15).Value = 0.85

# Row 4
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 30.748
$ws.Cells.Item(4, 4).Value = 61.476
$ws.Cells.Item(4, 5).Value = 2.85221546
$ws.Cells.Item(4, 6).Value = 0.17291004
$ws.Cells.Item(4, 7).Value = 0.08645518000000001
$ws.Cells.Item(4, 8).Value = 2.58611926
$ws.Cells.Item(4, 9).Value = 5.946698449130389
$ws.Cells.Item(4, 10).Value = 11.89505816789625
$ws.Cells.Item(4, 11).Value = 0.5518735913410133
$ws.Cells.Item(4, 12).Value = 0.03579685580843574
$ws.Cells.Item(4, 13).Value = 0.0178981127789971
$ws.Cells.Item(4, 14).Value = 0.3873361173245021
$ws.Cells.Item(4, 15).Value = 0.15

# Row 5
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = 48.81
$ws.Cells.Item(5, 4).Value = 95.57599999999999
$ws.Cells.Item(5, 5).Value = 1.81313864
$ws.Cells.Item(5, 6).Value = 0.11610202
$ws.Cells.Item(5, 7).Value = 0.05805088
$ws.Cells.Item(5, 8).Value = 2.77142636
$ws.Cells.Item(5, 9).Value = 8.241708045241959
$ws.Cells.Item(5, 10).Value = 14.91943037805436
$ws.Cells.Item(5, 11).Value = 0.2965792619205751
$ws.Cells.Item(5, 12).Value = 0.02113498621253413
$ws.Cells.Item(5, 13).Value = 0.0105675353243442
$ws.Cells.Item(5, 14).Value = 0.3576329599469703
$ws.Cells.Item(5, 15).Value = 0.85

# Row 6
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = 15.744
$ws.Cells.Item(6, 4).Value = 62.924
$ws.Cells.Item(6, 5).Value = 2.84737828
$ws.Cells.Item(6, 6).Value = 0.22580032
$ws.Cells.Item(6, 7).Value = 0.05645006
$ws.Cells.Item(6, 8).Value = 0.8657871399999999
$ws.Cells.Item(6, 9).Value = 3.781939996917
$ws.Cells.Item(6, 10).Value = 15.11380422515498
$ws.Cells.Item(6, 11).Value = 0.7046242549785309
$ws.Cells.Item(6, 12).Value = 0.06003703824657276
$ws.Cells.Item(6, 13).Value = 0.01500930852775579
$ws.Cells.Item(6, 14).Value = 0.2366686844066573
$ws.Cells.Item(6, 15).Value = 0.15

# Row 7
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = 25.348
$ws.Cells.Item(7, 4).Value = 93.38200000000001
$ws.Cells.Item(7, 5).Value = 1.85776186
$ws.Cells.Item(7, 6).Value = 0.15950234
$ws.Cells.Item(7, 7).Value = 0.03987556
$ws.Cells.Item(7, 8).Value = 0.98176622
$ws.Cells.Item(7, 9).Value = 5.114999750233599
$ws.Cells.Item(7, 10).Value = 15.04102859122887
$ws.Cells.Item(7, 11).Value = 0.3053325490125699
$ws.Cells.Item(7, 12).Value = 0.03645538373863952
$ws.Cells.Item(7, 13).Value = 0.009114106195886672
$ws.Cells.Item(7, 14).Value = 0.1766456612146574
$ws.Cells.Item(7, 15).Value = 0.85

# Row 8
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = 9.641999999999999
$ws.Cells.Item(8, 4).Value = 57.744
$ws.Cells.Item(8, 5).Value = 3.10425284
$ws.Cells.Item(8, 6).Value = 0.25772888
$ws.Cells.Item(8, 7).Value = 0.04295486
$ws.Cells.Item(8, 8).Value = 0.40571476
$ws.Cells.Item(8, 9).Value = 2.463732118535578
$ws.Cells.Item(8, 10).Value = 14.58904192497154
$ws.Cells.Item(8, 11).Value = 0.749643941158311
$ws.Cells.Item(8, 12).Value = 0.07542841647263143
$ws.Cells.Item(8, 13).Value = 0.01257127357187161
$ws.Cells.Item(8, 14).Value = 0.1366665586595862
$ws.Cells.Item(8, 15).Value = 0.15

# Row 9
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 17.61
$ws.Cells.Item(9, 4).Value = 88.34999999999999
$ws.Cells.Item(9, 5).Value = 1.97890148
$ws.Cells.Item(9, 6).Value = 0.17104938
$ws.Cells.Item(9, 7).Value = 0.0285082
$ws.Cells.Item(9, 8).Value = 0.48524196
$ws.Cells.Item(9, 9).Value = 4.664484411985886
$ws.Cells.Item(9, 10).Value = 16.32312435710624
$ws.Cells.Item(9, 11).Value = 0.3705008976022294
$ws.Cells.Item(9, 12).Value = 0.04396027155017505
$ws.Cells.Item(9, 13).Value = 0.007326587974131299
$ws.Cells.Item(9, 14).Value = 0.1230371086753129
$ws.Cells.Item(9, 15).Value = 0.85

# Row 10
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 7.274
$ws.Cells.Item(10, 4).Value = 58.048
$ws.Cells.Item(10, 5).Value = 3.0894319
$ws.Cells.Item(10, 6).Value = 0.25027196
$ws.Cells.Item(10, 7).Value = 0.0312839
$ws.Cells.Item(10, 8).Value = 0.2274484
$ws.Cells.Item(10, 9).Value = 1.802615556820528
$ws.Cells.Item(10, 10).Value = 14.33515571519045
$ws.Cells.Item(10, 11).Value = 0.7659593082350614
$ws.Cells.Item(10, 12).Value = 0.08033628923278562
$ws.Cells.Item(10, 13).Value = 0.01004224933635224
$ws.Cells.Item(10, 14).Value = 0.09301096028806183
$ws.Cells.Item(10, 15).Value = 0.15

# Row 11
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = 13.262
$ws.Cells.Item(11, 4).Value = 82.38
$ws.Cells.Item(11, 5).Value = 2.12173136
$ws.Cells.Item(11, 6).Value = 0.15903728
$ws.Cells.Item(11, 7).Value = 0.01987986
$ws.Cells.Item(11, 8).Value = 0.25712716
$ws.Cells.Item(11, 9).Value = 3.46603271270354
$ws.Cells.Item(11, 10).Value = 15.29429406222214
$ws.Cells.Item(11, 11).Value = 0.3900902000225816
$ws.Cells.Item(11, 12).Value = 0.04368976143113417
$ws.Cells.Item(11, 13).Value = 0.005461033267532637
$ws.Cells.Item(11, 14).Value = 0.07773990213978794
$ws.Cells.Item(11, 15).Value = 0.85

# Row 12
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = 5.992
$ws.Cells.Item(12, 4).Value = 59.748
$ws.Cells.Item(12, 5).Value = 3.0072619
$ws.Cells.Item(12, 6).Value = 0.22120214
$ws.Cells.Item(12, 7).Value = 0.02212039999999999
$ws.Cells.Item(12, 8).Value = 0.1337144
$ws.Cells.Item(12, 9).Value = 1.510131449241415
$ws.Cells.Item(12, 10).Value = 15.06567107333024
$ws.Cells.Item(12, 11).Value = 0.7547356655512579
$ws.Cells.Item(12, 12).Value = 0.07347401461515313
$ws.Cells.Item(12, 13).Value = 0.007347343025820841
$ws.Cells.Item(12, 14).Value = 0.06095312508111712
$ws.Cells.Item(12, 15).Value = 0.15

# Row 13
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = 11.086
$ws.Cells.Item(13, 4).Value = 75.96599999999999
$ws.Cells.Item(13, 5).Value = 2.30524572
$ws.Cells.Item(13, 6).Value = 0.14057886
$ws.Cells.Item(13, 7).Value = 0.01405778
$ws.Cells.Item(13, 8).Value = 0.15203048
$ws.Cells.Item(13, 9).Value = 3.081817515427673
$ws.Cells.Item(13, 10).Value = 14.61602203749026
$ws.Cells.Item(13, 11).Value = 0.435447527956801
$ws.Cells.Item(13, 12).Value = 0.04038502923700002
$ws.Cells.Item(13, 13).Value = 0.004038513495595653
$ws.Cells.Item(13, 14).Value = 0.05047845221523429
$ws.Cells.Item(13, 15).Value = 0.85
